# Append new ERA rows to the bottom of the data table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data to append: Index (col A), Name (col B), ERA (col C)
$newRows = @(
    @(871, "carl edwards", 2.76),
    @(872, "daniel lynch", 5.13),
    @(873, "duane underwood", 4.4),
    @(874, "frank german", 18),
    @(875, "j.t. chargois", 2.42),
    @(876, "jaime barría", 2.61),
    @(877, "lance mccullers", 2.27),
    @(878, "mark leiter", 3.99),
    @(879, "matt boyd", 1.35),
    @(880, "mike king", 2.29),
    @(881, "néstor cortés", 2.44),
    @(882, "nick martínez", 3.47),
    @(883, "travis lakins", 9.58),
    @(884, "vladimir gutiérrez", 7.61)
)

$startRow = 873
$lastExistingRow = 872

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $entry = $newRows[$i]

    $ws.Cells.Item($rowNum, 1).Value = $entry[0]
    $ws.Cells.Item($rowNum, 2).Value = $entry[1]
    $ws.Cells.Item($rowNum, 3).Value = $entry[2]

    # Match the styling used on the existing index column (column A) by
    # copying formatting from the last pre-existing row.
    $src = $ws.Cells.Item($lastExistingRow, 1)
    $dst = $ws.Cells.Item($rowNum, 1)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
